# Update cryptos list - price/volume refresh, and reorder ARBITRUM/ImmutableX and
# RenderToken/EnergySwap rows, as produced by the GitHub Actions data-refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.884.52'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.806.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.31'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4442'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3674'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07326'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8558'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.64'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.919.55'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.607'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.52'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.302'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07058'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008738'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.916.66'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.150'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.998'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.83'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.46'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.205'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.50'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08820'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("B31").Value = 'ARBITRUM'
$ws.Range("C31").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.175'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7490'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.933'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.459'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01965'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05192'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5327'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.860'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.012'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1687'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5152'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.433'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.59'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.985'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.44'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9993'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.665'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06322'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9176'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.60%  '
